# Update the dSF (column F) values per the re-pulled data / push all data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -1
    5  = 5
    6  = -7
    7  = -4
    8  = -2
    9  = 5
    10 = 1
    11 = -2
    13 = -4
    15 = 4
    16 = -5
    17 = 4
    18 = -7
    19 = -2
    20 = -2
    21 = -1
    22 = 1
    23 = 6
    24 = 1
    25 = 4
    26 = -2
    27 = 3
    29 = -2
    30 = -1
    31 = 1
    32 = -3
    33 = -2
    34 = -4
    35 = 1
    36 = -1
    37 = 1
    38 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
